# Generate Report for Handback
# Replace the first handback file's GUID/hash/timestamps and the second
# handback file's GUID across the Overview / zh-cn / de-de sheets, and
# refresh the hyperlinks that point at those filenames.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "859e9e9f-dcc1-4d34-9199-9ef345bb5a9b"
$newGuid1 = "885e7002-3dba-40ab-a7d3-33d242224785"
$oldGuid2 = "c70a462c-d133-44c2-9e86-4df3f6cc1309"
$newGuid2 = "ffff0f0c68c3-f6cb-482c-a2ac-55070b1e791d"

$newHash = "0271c48b8eceb70fe07976a53047e849de7936c1"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("B2").Value = "e2e\$newGuid1.md"
$ws.Range("G2").Value = "2016-08-16 11:03:35"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("B3").Value = "e2e\$newGuid2.md"
$ws.Range("G3").Value = "2016-08-16 11:03:35"

$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("B3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06a15762ae7f88aec57d2f2e3cc2bdb7f09ed630/e2e/$newGuid1.md", "", "", "e2e\$newGuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06a15762ae7f88aec57d2f2e3cc2bdb7f09ed630/e2e/$newGuid2.md", "", "", "e2e\$newGuid2.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-16 11:03:29"
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-16 11:03:57"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-16 11:03:29"
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-16 11:03:57"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("I2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()
$ws.Range("I3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06a15762ae7f88aec57d2f2e3cc2bdb7f09ed630/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/304054e4e8769161a5dfb9dc4cf1bac94868c5bf/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06a15762ae7f88aec57d2f2e3cc2bdb7f09ed630/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/304054e4e8769161a5dfb9dc4cf1bac94868c5bf/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("H2").Value = "2016-08-16 11:03:35"
$ws.Range("I2").Value = "$newGuid1.md"
$ws.Range("J2").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("K2").Value = "2016-08-16 11:04:12"

$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("H3").Value = "2016-08-16 11:03:35"
$ws.Range("I3").Value = "$newGuid2.md"
$ws.Range("J3").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("K3").Value = "2016-08-16 11:04:12"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("I2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()
$ws.Range("I3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06a15762ae7f88aec57d2f2e3cc2bdb7f09ed630/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8a4849aa1aa9ce3771c9fe6e0f620a14a1feb87e/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06a15762ae7f88aec57d2f2e3cc2bdb7f09ed630/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8a4849aa1aa9ce3771c9fe6e0f620a14a1feb87e/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null

Write-Output "done"
